$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("分区统计")

$ws.Range("B1:M1").ClearContents()
$ws.Range("B2:M3").ClearContents()
$ws.Range("B5:H12").ClearContents()

$ws.Activate()
$ws.Range("E24").Select() | Out-Null
